{"js": "// Replace the 25 \"three-digit number \u00d7 one-digit number\" expressions in the\n// practice-sheet table with their new values (problem text regenerated for a\n// different output, per the commit message / diff).\nconst replacements = [\n  [\"371\u00d72=\", \"517\u00d73=\"],\n  [\"743\u00d72=\", \"626\u00d78=\"],\n  [\"796\u00d74=\", \"648\u00d78=\"],\n  [\"131\u00d75=\", \"888\u00d79=\"],\n  [\"693\u00d72=\", \"197\u00d77=\"],\n  [\"398\u00d73=\", \"333\u00d76=\"],\n  [\"429\u00d74=\", \"803\u00d76=\"],\n  [\"925\u00d75=\", \"722\u00d75=\"],\n  [\"622\u00d77=\", \"779\u00d73=\"],\n  [\"826\u00d79=\", \"498\u00d75=\"],\n  [\"838\u00d76=\", \"161\u00d75=\"],\n  [\"455\u00d77=\", \"680\u00d78=\"],\n  [\"171\u00d74=\", \"612\u00d78=\"],\n  [\"809\u00d72=\", \"237\u00d75=\"],\n  [\"514\u00d76=\", \"356\u00d76=\"],\n  [\"210\u00d72=\", \"252\u00d77=\"],\n  [\"712\u00d73=\", \"164\u00d74=\"],\n  [\"358\u00d79=\", \"503\u00d74=\"],\n  [\"494\u00d73=\", \"353\u00d79=\"],\n  [\"474\u00d75=\", \"884\u00d73=\"],\n  [\"925\u00d74=\", \"233\u00d76=\"],\n  [\"774\u00d78=\", \"728\u00d73=\"],\n  [\"140\u00d79=\", \"125\u00d73=\"],\n  [\"595\u00d74=\", \"782\u00d74=\"],\n  [\"984\u00d76=\", \"631\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit number x one-digit number\" expressions in the\n# practice-sheet table with their new values (problem text regenerated for a\n# different output, per the commit message / diff).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"371\u00d72=\", \"517\u00d73=\"),\n    @(\"743\u00d72=\", \"626\u00d78=\"),\n    @(\"796\u00d74=\", \"648\u00d78=\"),\n    @(\"131\u00d75=\", \"888\u00d79=\"),\n    @(\"693\u00d72=\", \"197\u00d77=\"),\n    @(\"398\u00d73=\", \"333\u00d76=\"),\n    @(\"429\u00d74=\", \"803\u00d76=\"),\n    @(\"925\u00d75=\", \"722\u00d75=\"),\n    @(\"622\u00d77=\", \"779\u00d73=\"),\n    @(\"826\u00d79=\", \"498\u00d75=\"),\n    @(\"838\u00d76=\", \"161\u00d75=\"),\n    @(\"455\u00d77=\", \"680\u00d78=\"),\n    @(\"171\u00d74=\", \"612\u00d78=\"),\n    @(\"809\u00d72=\", \"237\u00d75=\"),\n    @(\"514\u00d76=\", \"356\u00d76=\"),\n    @(\"210\u00d72=\", \"252\u00d77=\"),\n    @(\"712\u00d73=\", \"164\u00d74=\"),\n    @(\"358\u00d79=\", \"503\u00d74=\"),\n    @(\"494\u00d73=\", \"353\u00d79=\"),\n    @(\"474\u00d75=\", \"884\u00d73=\"),\n    @(\"925\u00d74=\", \"233\u00d76=\"),\n    @(\"774\u00d78=\", \"728\u00d73=\"),\n    @(\"140\u00d79=\", \"125\u00d73=\"),\n    @(\"595\u00d74=\", \"782\u00d74=\"),\n    @(\"984\u00d76=\", \"631\u00d73=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n}\n"}
